$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''57.868.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.85%  '

# Row 3
$ws.Range("D3").Value = '''3.046.55'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.51%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = '''517.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.16%  '

# Row 6
$ws.Range("D6").Value = '''140.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.70%  '

# Row 7
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("E8").Value = '  +3.69%  '

# Row 9
$ws.Range("D9").Value = '''7.53'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.82%  '

# Row 10
$ws.Range("E10").Value = '  +4.50%  '

# Row 11
$ws.Range("E11").Value = '  +5.49%  '

# Row 12
$ws.Range("D12").Value = '''3.571.64'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.55%  '

# Row 13
$ws.Range("E13").Value = '  +2.29%  '

# Row 14
$ws.Range("D14").Value = '''26.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.42%  '

# Row 15
$ws.Range("E15").Value = '  +11.54%  '

# Row 16
$ws.Range("D16").Value = '''57.855.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.79%  '

# Row 17
$ws.Range("D17").Value = '''6.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +10.26%  '

# Row 18
$ws.Range("D18").Value = '''3.038.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.08%  '

# Row 19
$ws.Range("D19").Value = '''13.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.67%  '

# Row 20
$ws.Range("D20").Value = '''8.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.11%  '

# Row 21
$ws.Range("D21").Value = '''337.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.39%  '

# Row 22
$ws.Range("E22").Value = '  +1.48%  '

# Row 23
$ws.Range("D23").Value = '''1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '

# Row 24
$ws.Range("D24").Value = '''0.502'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.88%  '

# Row 25
$ws.Range("D25").Value = '''64.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.92%  '

# Row 26
$ws.Range("E26").Value = '  +4.23%  '

# Row 27
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '''0.0₃0948'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.15%  '

# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '''0.990'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.03%  '

# Row 29
$ws.Range("D29").Value = '''6.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.17%  '

# Row 30
$ws.Range("D30").Value = '''7.52'
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = '  +4.55%  '

# Row 32
$ws.Range("E32").Value = '  +3.23%  '

# Row 33
$ws.Range("D33").Value = '''20.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.86%  '

# Row 34
$ws.Range("D34").Value = '''156.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.07%  '

# Row 35
$ws.Range("D35").Value = '''4.76'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.84%  '

# Row 36
$ws.Range("D36").Value = '''5.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.85%  '

# Row 37
$ws.Range("E37").Value = '  +2.50%  '

# Row 38
$ws.Range("D38").Value = '''25.11'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.47%  '

# Row 39
$ws.Range("D39").Value = '''0.0692'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.78%  '

# Row 40
$ws.Range("D40").Value = '''3.082.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.49%  '

# Row 41
$ws.Range("D41").Value = '''37.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.34%  '

# Row 42
$ws.Range("D42").Value = '''3.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.07%  '

# Row 43
$ws.Range("E43").Value = '  +0.00%  '

# Row 44
$ws.Range("E44").Value = '  +3.67%  '

# Row 45
$ws.Range("D45").Value = '''2.312.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.78%  '

# Row 46
$ws.Range("E46").Value = '  +3.78%  '

# Row 47
$ws.Range("E47").Value = '  +2.03%  '

# Row 48
$ws.Range("D48").Value = '''6.08'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.64%  '

# Row 49
$ws.Range("E49").Value = '  +3.16%  '

# Row 50
$ws.Range("D50").Value = '''19.70'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.01%  '

# Row 51
$ws.Range("E51").Value = '  -4.51%  '
